# Update cryptocurrency price/volume data on the "cryptos" sheet.
# Commit message: Updated cryptos list on Fri Mar 24 19:22:30 UTC 2023 with GitHub Actions
#
# All the "Price" (column D) and "Volume(1h)" (column E) cells in this sheet
# are plain text (inline strings) in the source workbook, not real numbers -
# many of them look like numbers ("0.9995") or even use a dotted
# thousands-style grouping ("27.787.17"). Excel's Range.Value setter will
# silently coerce a numeric-looking string into a real number (dropping
# trailing zeros, switching to scientific notation for tiny values, etc.),
# so column D is temporarily forced to Text format while the new values are
# written, then restored to the default "Normal" style so the saved file
# doesn't carry a stray text-format style around.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text so numeric-looking strings are kept verbatim.
$ws.Range("D2:D51").NumberFormat = "@"

# ---- Rows 2-36: price (D) and volume-1h (E) refresh ----
$ws.Range("D2").Value = "27.816.27"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "1.761.33"
$ws.Range("E3").Value = "  -2.60%  "
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").Value = "321.72"
$ws.Range("E5").Value = "  -2.44%  "
$ws.Range("D6").Value = "0.9988"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D7").Value = "0.4240"
$ws.Range("E7").Value = "  -4.04%  "
$ws.Range("D8").Value = "0.3630"
$ws.Range("E8").Value = "  -2.48%  "
$ws.Range("E9").Value = "  -5.15%  "
$ws.Range("D10").Value = "0.07473"
$ws.Range("E10").Value = "  -2.84%  "
$ws.Range("D11").Value = "1.087"
$ws.Range("E11").Value = "  -2.82%  "
$ws.Range("D12").Value = "0.9998"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").Value = "20.72"
$ws.Range("E13").Value = "  -5.37%  "
$ws.Range("D14").Value = "6.070"
$ws.Range("E14").Value = "  -3.47%  "
$ws.Range("D15").Value = "7.294"
$ws.Range("E15").Value = "  -2.48%  "
$ws.Range("D16").Value = "1.753.38"
$ws.Range("E16").Value = "  -4.12%  "
$ws.Range("E17").Value = "  -2.58%  "
$ws.Range("D18").Value = "0.00001055"
$ws.Range("E18").Value = "  -2.51%  "
$ws.Range("D19").Value = "0.06369"
$ws.Range("E19").Value = "  -1.78%  "
$ws.Range("D20").Value = "0.9996"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "17.03"
$ws.Range("E21").Value = "  -2.41%  "
$ws.Range("D22").Value = "5.933"
$ws.Range("E22").Value = "  -5.17%  "
$ws.Range("D23").Value = "27.828.19"
$ws.Range("E23").Value = "  -0.99%  "
$ws.Range("D24").Value = "11.22"
$ws.Range("E24").Value = "  -3.77%  "
$ws.Range("D25").Value = "2.101"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "157.48"
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("D27").Value = "20.21"
$ws.Range("E27").Value = "  -1.85%  "
$ws.Range("D28").Value = "1.958.65"
$ws.Range("E28").Value = "  -3.39%  "
$ws.Range("D29").Value = "2.135"
$ws.Range("E29").Value = "  -8.02%  "
$ws.Range("D30").Value = "124.03"
$ws.Range("E30").Value = "  -2.36%  "
$ws.Range("D31").Value = "1.117"
$ws.Range("E31").Value = "  -6.68%  "
$ws.Range("D32").Value = "3.684"
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").Value = "5.562"
$ws.Range("E33").Value = "  -4.64%  "
$ws.Range("D34").Value = "0.08853"
$ws.Range("E34").Value = "  -4.08%  "
$ws.Range("D35").Value = "12.23"
$ws.Range("E35").Value = "  -6.16%  "
$ws.Range("D36").Value = "0.02293"
$ws.Range("E36").Value = "  -1.88%  "

# ---- Rows 37-40: ranking reshuffle ----
# Hedera and Algorand swap positions (37 <-> 38), and TheSandbox /
# InternetComputer(DFINITY) swap positions (39 <-> 40), with refreshed values.
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "0.2100"
$ws.Range("E37").Value = "  -3.06%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.06043"
$ws.Range("E38").Value = "  -2.37%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "4.970"
$ws.Range("E39").Value = "  -3.74%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.6326"
$ws.Range("E40").Value = "  -3.49%  "

# ---- Rows 41-51: price (D) and volume-1h (E) refresh ----
$ws.Range("D41").Value = "1.175"
$ws.Range("E41").Value = "  -1.57%  "
$ws.Range("D42").Value = "0.9986"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("D43").Value = "7.874"
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("D44").Value = "1.396"
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("E45").Value = "  -5.33%  "
$ws.Range("D46").Value = "0.5866"
$ws.Range("E46").Value = "  -3.34%  "
$ws.Range("D47").Value = "3.684"
$ws.Range("E47").Value = "  -2.16%  "
$ws.Range("D48").Value = "1.988"
$ws.Range("E48").Value = "  -2.22%  "
$ws.Range("D49").Value = "122.96"
$ws.Range("E49").Value = "  -2.98%  "
$ws.Range("D50").Value = "1.180"
$ws.Range("E50").Value = "  +2.64%  "
$ws.Range("D51").Value = "0.06822"
$ws.Range("E51").Value = "  -2.12%  "

# Restore the default (unstyled) look for column D so the saved cells don't
# carry an explicit Text-format style around (matches the original file,
# where these cells have no style index at all).
$ws.Range("D2:D51").Style = "Normal"

Write-Output "cryptos sheet updated"
